$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "A19" (physical sheet20.xml) - add three rows of TxHash values
# in column A (column B already has the matching ChainID values).
# Order matters: these three new strings must land at shared-string
# indices 61, 62, 63 (in that order).
# ---------------------------------------------------------------------
$wsA19 = $wb.Worksheets.Item("A19")

# Copy the existing "TxHash" cell formatting (style s=2) down into the
# new A5:A7 cells before filling in the values.
$wsA19.Range("A2:A4").Copy()
$wsA19.Range("A5:A7").PasteSpecial(-4122)

$wsA19.Range("A5").Value = "2301B6D3CA98F08D36B9D8A76F87032872AC9D45DF8A4FE5EC266ACE80277D1B"
$wsA19.Range("A6").Value = "BB4D299DB10FA4AD9CC1E49D81A0FC9759CA0A9A4D15E5E68CC4DB17C4E7D62C"
$wsA19.Range("A7").Value = "0F2CAFAF63981422EC83854FE0FD1493EBF2CC96FAC5E3D8F9FFF06ED4831CB5"

# ---------------------------------------------------------------------
# Sheet "A20" (physical sheet21.xml) - replace the placeholder header
# text in column A with real TxHash values, fill in column B ChainID
# values, and append three more evidence rows.
# New shared strings introduced here must append in this order:
#   "gon-flixnet-1" (64), then the six hashes below (65..70)
# ---------------------------------------------------------------------
$wsA20 = $wb.Worksheets.Item("A20")

# Give the new B3:B7 cells the same formatting (style s=5, left aligned)
# as the existing B2 cell had before its style also moves to s=5.
$wsA19.Range("B2").Copy()
$wsA20.Range("B2:B7").PasteSpecial(-4122)

# Extend column A's formatting (style s=2) down to the two new rows.
$wsA20.Range("A4").Copy()
$wsA20.Range("A5:A7").PasteSpecial(-4122)

# Column B (ChainID) values - row 5's "gon-flixnet-1" is a brand-new
# shared string and must be written before the column-A hashes below.
$wsA20.Range("B5").Value = "gon-flixnet-1"
$wsA20.Range("B2").Value = "gon-irishub-1"
$wsA20.Range("B3").Value = "uptick_7000-2"
$wsA20.Range("B4").Value = "elgafar-1"
$wsA20.Range("B6").Value = "elgafar-1"
$wsA20.Range("B7").Value = "uptick_7000-2"

# Column A (TxHash) values, in row order - these six strings are new
# and must append to the shared-string table in this order (65..70).
$wsA20.Range("A2").Value = "6B8C02B83AA26A500F8624C64283A13E92A3D17D1C3D57F9A3AB0E55315E8FBC"
$wsA20.Range("A3").Value = "C1083554FDB1E2CF6B21E9389F310BCD56139F6CCCF0B82ADA139790ACDFD9F8"
$wsA20.Range("A4").Value = "37478751EC7DBB70AEFFDE6011FBAF4C7307122B45E495045C06DA9EE7376959"
$wsA20.Range("A5").Value = "CA1ECC7C0951461718BE2DD10254CB9A4965C79579C66766DA1191E5FB93A1FE"
$wsA20.Range("A6").Value = "3C26C8F553106EF7162F133067CFB230D688509EE2DA75E78D2B39FCFEFD6570"
$wsA20.Range("A7").Value = "BBF063AFA163396AE2EAB548220C6E0CEC77D24FA98B780DFE62837FFE54DC22"

# ---------------------------------------------------------------------
# Window / tab-selection state: the active sheet moves from A19 to A20,
# and the remembered selection cell on each sheet shifts down by a row.
# ---------------------------------------------------------------------
$wsA19.Range("C9").Select()
$wsA20.Activate()
$wsA20.Range("C4").Select()
